# Apply the changes described by the diff:
# 1. Metadata sheet: update the "Date" value in B8.
# 2. Elements sheet: swap the two "Mapping" columns (AK <-> AL), including
#    their header text, their per-row values, and their column widths.

$wb = $excel.ActiveWorkbook

# --- 1. Update Date on Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Swap Mapping columns on Elements sheet ---
$els = $wb.Worksheets.Item("Elements")

# Swap cell contents/values for rows 1 (header) through 6 (last data row)
for ($r = 1; $r -le 6; $r++) {
    $akCell = $els.Range("AK$r")
    $alCell = $els.Range("AL$r")

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the column widths of columns AK (37) and AL (38).
# (Original widths: AK = 24.98046875, AL = 74.24609375 -- taken from the
# workbook's stored column definitions, since the COM ColumnWidth getter
# does not reliably reflect custom widths loaded from the file.)
$els.Columns.Item(37).ColumnWidth = 74.24609375
$els.Columns.Item(38).ColumnWidth = 24.98046875
